$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Fix a typo in the "Harapan Orang Tua" column for row 13 (student #12):
# "Mendapatkan ilmu yang bermanfa.at " -> "Mendapatkan ilmu yang bermanfaat "
$ws.Range("Z13").Value = "Mendapatkan ilmu yang bermanfaat "

# Rename the photo filenames in column AE (Foto) from the original
# upload names to sequential numbers: row 2 -> "1.jpg", row 3 -> "2.jpg", ...
for ($row = 2; $row -le 33; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 31).Value = "$n.jpg"
}

# Update the view: move the selection and zoom level to match the saved view
$ws.Range("AI7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 115
